$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16
$ws.Range("D16").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E16").Value = "['Normal']"

# Row 20
$ws.Range("D20").Value = "[0, 0, 0, 0, 0, 0, 1]"
$ws.Range("E20").Value = "['SoftwareFault']"

# Row 29
$ws.Range("D29").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E29").Value = "['Normal']"

# Row 97
$ws.Range("D97").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E97").Value = "['Normal']"

# Row 113
$ws.Range("D113").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E113").Value = "['Normal', 'HardwareFault']"
